$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 originally carried a bespoke bold/no-alignment row-level style
# (s="4" customFormat="1" on the row, s="3" on each cell). The new layout
# drops that special casing so row 2 matches the other data rows (s="1",
# plain centered alignment, no per-row custom format).
$ws.Rows(2).ClearFormats()
$ws.Range("A2:G2").HorizontalAlignment = -4108

# New header cells for the two similarity columns
$ws.Range("H1").Value = "op sim"
$ws.Range("I1").Value = "data sim"
$ws.Range("H1:I1").Font.Bold = $true
$ws.Range("H1:I1").HorizontalAlignment = -4108

# New "op sim" / "data sim" data for rows 2-11, matching the rest of the
# table's plain centered style
$ws.Range("H2:I11").HorizontalAlignment = -4108

$hvals = @(-0.0265, 0.6913, 0.786, 0.8505, 0.8453, 0.8733, 0.9491, 0.9329, 0.9253, 0.9436)
$ivals = @(0.1784, 0.687, 0.7602, 0.7967, 0.8109, 0.8392, 0.8924, 0.8924, 0.8924, 0.9027)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $hvals[$i]
    $ws.Cells.Item($row, 9).Value = $ivals[$i]
}

# Move the active selection to reflect where the author ended up (H15)
$ws.Range("H15").Select()
